# 2025 author list / membership refresh for the Imperial collaboration list.
# - Updates two existing members' email addresses to include extra addresses
#   (keeping JISCMail / collaboration DB in sync).
# - Appends six new collaborators (rows 24-29), each with a mailto:
#   hyperlink on the email-address column (matching the convention already
#   used elsewhere in the sheet).
# - Grows Table1 to cover the new rows and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Refresh email address for Pat Price first (matches original authoring
#    order), then fill in the first new collaborator row, then refresh
#    Kenneth Long's email.
# ---------------------------------------------------------------------------
$ws.Range("E22").Value = "p.price@imperial.ac.uk; patprice@patprice.co.uk"

$ws.Range("A24").Value = "Mx"
$ws.Range("B24").Value = "Suriu"
$ws.Range("C24").Value = "Liy"
$ws.Range("D24").Value = "S."
$ws.Range("E24").Value = "s.lyu23@imperial.ac.uk "
$ws.Range("F24").Value = "S.Lyu"
$ws.Range("G24").Value = "Imperial BioEng"
$ws.Range("H24").Value = "Department of Bioengineering, Imperial College London, Exhibition Road, London, SW7 2AZ, UK"
$ws.Range("I24").Value = 0

$ws.Range("E19").Value = "k.long@imperial.ac.uk; ken.long@STFC.ac.uk; longkr@ic.ac.uk"

# ---------------------------------------------------------------------------
# 2. Remaining new collaborators (rows 25-29).
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = "C."
$ws.Range("G25").Value = "Imperial-Physics"
$ws.Range("H25").Value = "Department of Physics, Imperial College London, Exhibition Road, London, SW7 2AZ, UK"
$ws.Range("A25").Value = "Mr."
$ws.Range("B25").Value = "Calvin"
$ws.Range("C25").Value = "Dyson"
$ws.Range("E25").Value = "c.dyson24@imperial.ac.uk"
$ws.Range("F25").Value = "C.Dyson"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = "JohnAdamsInstitute-Imperial"
$ws.Range("K25").Value = "The John Adams Institute for Accelerator Science, Department of Physics, Imperial College London, Exhibition Road, London SW7 2AZ, UK"

$ws.Range("G26").Value = "Imperial-Physics"
$ws.Range("H26").Value = "Department of Physics, Imperial College London, Exhibition Road, London, SW7 2AZ, UK"
$ws.Range("A26").Value = "Dr."
$ws.Range("B26").Value = "Titus"
$ws.Range("C26").Value = "Dascalu"
$ws.Range("D26").Value = "T.S."
$ws.Range("E26").Value = "t.dascalu19@imperial.ac.uk"
$ws.Range("F26").Value = "T.S.Dascalu"
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = "JohnAdamsInstitute-Imperial"
$ws.Range("K26").Value = "The John Adams Institute for Accelerator Science, Department of Physics, Imperial College London, Exhibition Road, London SW7 2AZ, UK"

$ws.Range("G27").Value = "Imperial BioEng"
$ws.Range("H27").Value = "Department of Bioengineering, Imperial College London, Exhibition Road, London, SW7 2AZ, UK"
$ws.Range("A27").Value = "Mx."
$ws.Range("B27").Value = "Xiangyi"
$ws.Range("C27").Value = "Chen"
$ws.Range("D27").Value = "X."
$ws.Range("E27").Value = "x.chen24@imperial.ac.uk"
$ws.Range("F27").Value = "X.Chen"
$ws.Range("I27").Value = 0

$ws.Range("G28").Value = "Imperial-Physics"
$ws.Range("H28").Value = "Department of Physics, Imperial College London, Exhibition Road, London, SW7 2AZ, UK"
$ws.Range("A28").Value = "Mx."
$ws.Range("B28").Value = "Ginevra"
$ws.Range("E28").Value = "ginevra.casati18@imperial.ac.uk"
$ws.Range("C28").Value = "Casati"
$ws.Range("D28").Value = "G."
$ws.Range("F28").Value = "G.Casati"
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = "JohnAdamsInstitute-Imperial"
$ws.Range("K28").Value = "The John Adams Institute for Accelerator Science, Department of Physics, Imperial College London, Exhibition Road, London SW7 2AZ, UK"

$ws.Range("D29").Value = "D."
$ws.Range("G29").Value = "Imperial-SurgCanc"
$ws.Range("H29").Value = "Department of Surgery and Cancer, Imperial College, Hammersmith Hospital London, W12 0NN, UK"
$ws.Range("A29").Value = "Mr."
$ws.Range("B29").Value = "Diaza"
$ws.Range("E29").Value = "d.ariyanto24@imperial.ac.uk"
$ws.Range("F29").Value = "D.Ariyanto"
$ws.Range("C29").Value = "Ariyanto"
$ws.Range("I29").Value = 0

# ---------------------------------------------------------------------------
# 3. Hyperlink the email-address cell for each new row to a mailto: link,
#    same pattern already used elsewhere in the sheet for collaborator
#    emails.
# ---------------------------------------------------------------------------
$emailCells = @("E24", "E25", "E26", "E27", "E28", "E29")
$emailAddrs = @(
    "s.lyu23@imperial.ac.uk",
    "c.dyson24@imperial.ac.uk",
    "t.dascalu19@imperial.ac.uk",
    "x.chen24@imperial.ac.uk",
    "ginevra.casati18@imperial.ac.uk",
    "d.ariyanto24@imperial.ac.uk"
)
for ($i = 0; $i -lt $emailCells.Length; $i++) {
    $ws.Hyperlinks.Add($ws.Range($emailCells[$i]), "mailto:" + $emailAddrs[$i])
}

# ---------------------------------------------------------------------------
# 4. Grow Table1 to include the newly added rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P29"))

# ---------------------------------------------------------------------------
# 5. Move the active selection (matches the author's saved cursor position).
# ---------------------------------------------------------------------------
$ws.Range("I8").Select()
